$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.755.12'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.650.93'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.00'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '1.880.77'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +3.40%  '
$ws.Range('D14').Value = '1.649.69'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('E16').Value = '  +5.55%  '
$ws.Range('D17').Value = '26.815.20'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '0.0₃0757'
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.43'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.42'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.37'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +12.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.07'
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.121'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('E28').Value = '  +4.16%  '
$ws.Range('E29').Value = '  +3.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0522'
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.44'
$ws.Range('E32').Value = '  +4.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.06'
$ws.Range('E33').Value = '  +4.12%  '
$ws.Range('E34').Value = '  +3.87%  '
$ws.Range('D35').Value = '1.294.26'
$ws.Range('E35').Value = '  +8.27%  '
$ws.Range('E36').Value = '  +5.20%  '
$ws.Range('E37').Value = '  +1.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.837'
$ws.Range('E38').Value = '  +3.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.529'
$ws.Range('E39').Value = '  +4.77%  '
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.818'
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('E42').Value = '  -2.73%  '
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('D44').Value = '1.792.40'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.89'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.80'
$ws.Range('E46').Value = '  +9.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.62'
$ws.Range('E47').Value = '  +5.41%  '
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.80'
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0981'
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.409'
$ws.Range('E51').Value = '  -0.37%  '
